$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.034.24'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.553.40'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.19'
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3942'
$ws.Range("E7").Value = '  +3.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3222'
$ws.Range("E8").Value = '  -2.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.85'
$ws.Range("E9").Value = '  -1.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07238'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.074'
$ws.Range("E11").Value = '  -5.97%  '
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.662'
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.74'
$ws.Range("E14").Value = '  -7.30%  '
$ws.Range("E15").Value = '  +5.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.619'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").Value = '1.556.30'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06570'
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.39'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.267'
$ws.Range("E21").Value = '  -2.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.50'
$ws.Range("E22").Value = '  -4.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.29'
$ws.Range("E23").Value = '  -3.69%  '
$ws.Range("D24").Value = '22.047.25'
$ws.Range("E24").Value = '  -0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.368'
$ws.Range("E25").Value = '  +4.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.407'
$ws.Range("E26").Value = '  -6.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.99'
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.51'
$ws.Range("E28").Value = '  -4.11%  '
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("D30").Value = '1.726.21'
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.59'
$ws.Range("E31").Value = '  -3.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9735'
$ws.Range("E32").Value = '  -10.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.821'
$ws.Range("E33").Value = '  -1.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08339'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.068'
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("E36").Value = '  -16.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02259'
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("E38").Value = '  -4.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06005'
$ws.Range("E39").Value = '  -4.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.209'
$ws.Range("E40").Value = '  -2.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2033'
$ws.Range("E41").Value = '  -5.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.65'
$ws.Range("E43").Value = '  -3.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5799'
$ws.Range("E44").Value = '  -4.39%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.04'
$ws.Range("E45").Value = '  -5.74%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.743'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5563'
$ws.Range("E47").Value = '  -5.37%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.895'
$ws.Range("E48").Value = '  -3.74%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.65'
$ws.Range("E49").Value = '  -4.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.133'
$ws.Range("E50").Value = '  -3.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06819'
$ws.Range("E51").Value = '  -3.44%  '
